# "Restored from revision of admin ... Type: SAVE" - reverts cell C10 on
# the "Rules" sheet from 18 back to 1 (a stray earlier edit undone).

$wb = $excel.ActiveWorkbook

$ws = $null
try {
    $ws = $wb.Sheets("Rules")
} catch {
    $ws = $null
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Range("C10").Value = 1
